{"js": "// Remove the \"Appendix: Quick prototype\" section: its Heading 2 title,\n// the blank paragraph after it, the \"Figure: PDF page 1\" caption, and the\n// paragraph holding the embedded screenshot (and therefore the image\n// itself, plus the stray blank line right after the image).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"Appendix: Quick prototype\" Heading 2 paragraph.\nlet startIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const p = paras.items[i];\n  if (p.style === \"Heading 2\" && p.text.trim() === \"Appendix: Quick prototype\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error('Could not find \"Appendix: Quick prototype\" heading.');\n}\n\n// Find the next Heading 2 (the following \"Appendix: Links\" section) so we\n// know where the prototype appendix block ends; delete everything from the\n// \"Appendix: Quick prototype\" heading up to (but not including) that\n// paragraph.\nlet endIndex = paras.items.length;\nfor (let i = startIndex + 1; i < paras.items.length; i++) {\n  if (paras.items[i].style === \"Heading 2\") {\n    endIndex = i;\n    break;\n  }\n}\n\nfor (let i = endIndex - 1; i >= startIndex; i--) {\n  paras.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Appendix: Quick prototype\" section: its Heading 2 title,\n# the blank paragraph after it, the \"Figure: PDF page 1\" caption, and the\n# paragraph holding the embedded screenshot (and therefore the image\n# itself, plus the stray blank line right after the image).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$startIndex = -1\n$endIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    $txt = $p.Range.Text.Trim()\n    if ($styleName -eq \"Heading 2\" -and $txt -eq \"Appendix: Quick prototype\") {\n        $startIndex = $i\n    }\n    elseif ($startIndex -ge 1 -and $endIndex -eq -1 -and $styleName -eq \"Heading 2\") {\n        $endIndex = $i\n        break\n    }\n}\n\nif ($startIndex -eq -1) {\n    throw \"Could not find 'Appendix: Quick prototype' heading.\"\n}\n\n$startPara = $d.Paragraphs.Item($startIndex)\n\nif ($endIndex -eq -1) {\n    # No following Heading 2 found; delete through the end of the document.\n    $rng = $d.Range($startPara.Range.Start, $d.Content.End)\n}\nelse {\n    $endPara = $d.Paragraphs.Item($endIndex)\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.Start)\n}\n\n$rng.Delete()\n"}
